# This edit swaps the data contents of row 3 and row 5 on the active sheet
# (two species-observation records that exchanged their row positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Save the "before" values of row 3 that are needed to populate row 5 ---
$A3 = $ws.Range("A3").Value2
$B3 = $ws.Range("B3").Value2
$D3 = $ws.Range("D3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2
$P3 = $ws.Range("P3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2
$S3 = $ws.Range("S3").Value2
$AW3 = $ws.Range("AW3").Value2
$AX3 = $ws.Range("AX3").Value2

# --- Save the "before" values of row 5 that are needed to populate row 3 ---
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$D5 = $ws.Range("D5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$P5 = $ws.Range("P5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2
$S5 = $ws.Range("S5").Value2
$AW5 = $ws.Range("AW5").Value2
$AX5 = $ws.Range("AX5").Value2

# --- Write row 5's old data into row 3 ---
$ws.Range("A3").Value = $A5
$ws.Range("B3").Value = $B5
$ws.Range("D3").Value = $D5
$ws.Range("E3").Value = $E5
$ws.Range("F3").Value = $F5
$ws.Range("G3").Value = $G5
$ws.Range("H3").Value = $H5
$ws.Range("P3").Value = $P5
$ws.Range("Q3").Value = $Q5
$ws.Range("R3").Value = $R5
$ws.Range("S3").Value = $S5
$ws.Range("AW3").Value = $AW5
$ws.Range("AX3").Value = $AX5

# --- Write row 3's old data into row 5 ---
$ws.Range("A5").Value = $A3
$ws.Range("B5").Value = $B3
$ws.Range("D5").Value = $D3
$ws.Range("E5").Value = $E3
$ws.Range("F5").Value = $F3
$ws.Range("G5").Value = $G3
$ws.Range("H5").Value = $H3
$ws.Range("P5").Value = $P3
$ws.Range("Q5").Value = $Q3
$ws.Range("R5").Value = $R3
$ws.Range("S5").Value = $S3
$ws.Range("AW5").Value = $AW3
$ws.Range("AX5").Value = $AX3

# --- Row 3 no longer carries these (previously empty) fields; row 5 now does ---
# (the source records differ slightly in which blank columns were exported)
$ws.Range("J3").Clear()
$ws.Range("L3").Clear()
$ws.Range("N3").Clear()
$ws.Range("AF3").Clear()

$ws.Range("J5").Style = "Normal"
$ws.Range("L5").Style = "Normal"
$ws.Range("N5").Style = "Normal"
$ws.Range("AF5").Style = "Normal"
